$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Supplier PN for row 9 (10u / C402 capacitor)
$ws.Range("G9").Value = "81-GRM188R60J106ME47"

# Update Manufacturer PN and Pricing for row 11 (4.7u / C203 capacitor)
$ws.Range("I11").Value = "GRM188R60J106ME47D"
$ws.Range("J11").Value = "0.163/0.096/--"

# Scroll the sheet view so column H is the first visible column
$ws.Application.ActiveWindow.ScrollColumn = 8
